# "Adicionados valores para 641" — fill in the B/C measurement columns for
# the 641 series (rows 7-14), clear the old G7/H7 pair, add the new
# "0.357393 " text reading in C10, and move the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: replace the previous placeholder pair with the real reading ---
$ws.Range("B7").Value = 0.254189
$ws.Range("C7").Value = 2.117

# G7/H7 no longer hold a value for this series
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()

# --- Rows 8-14: fill in the newly measured values ---
$ws.Range("B8").Value = 0.129997
$ws.Range("C8").Value = 1.074014

$ws.Range("B9").Value = 0.065333
$ws.Range("C9").Value = 0.558948

$ws.Range("B10").Value = 0.036775
# C10 is a text reading ("0.357393 ", trailing space) rather than a number.
# Build it as a text formula, then paste-special the computed value back on
# top of itself so it is stored as a plain shared string (not a live
# formula) while keeping the cell's existing style untouched.
$ws.Range("C10").Formula = '="0.357393 "'
$ws.Range("C10").Copy()
$ws.Range("C10").PasteSpecial(-4163) | Out-Null

$ws.Range("B11").Value = 0.031687
$ws.Range("C11").Value = 0.301891

$ws.Range("B12").Value = 0.027754
$ws.Range("C12").Value = 0.241983

$ws.Range("B13").Value = 0.061819
$ws.Range("C13").Value = 0.261352

$ws.Range("B14").Value = 0.492241
$ws.Range("C14").Value = 1.12211

# --- View: scroll the sheet down a bit and land the selection on E14 ---
[void]($excel.ActiveWindow.ScrollRow = 4)
$ws.Range("E14").Select() | Out-Null
